$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some new values in column D (e.g. '6.17') look numeric to Excel's
# auto-detection, but must stay plain text to match the original
# inline-string cell type. Temporarily force a text number format on
# those specific cells, assign the values, then restore each cell's
# original style so the saved style is unchanged.
$textCells = @("D5", "D6", "D11", "D12", "D14", "D17", "D19", "D20", "D21", "D22", "D23", "D25", "D27", "D28", "D30", "D32", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D42", "D44", "D45", "D46", "D51")
$origStyles = @{}
foreach ($addr in $textCells) {
    $cell = $ws.Range($addr)
    $origStyles[$addr] = $cell.Style
    $cell.NumberFormat = "@"
}

$ws.Range("D2").Value = '64.677.47'
$ws.Range("E2").Value = '  +0.45%  '
$ws.Range("D3").Value = '3.152.83'
$ws.Range("E3").Value = '  +1.77%  '
$ws.Range("D5").Value = '571.28'
$ws.Range("E5").Value = '  +2.07%  '
$ws.Range("D6").Value = '150.61'
$ws.Range("E6").Value = '  +4.15%  '
$ws.Range("D8").Value = '3.150.55'
$ws.Range("E8").Value = '  +1.77%  '
$ws.Range("E9").Value = '  +4.51%  '
$ws.Range("E10").Value = '  +5.21%  '
$ws.Range("D11").Value = '6.17'
$ws.Range("E11").Value = '  +0.67%  '
$ws.Range("D12").Value = '0.504'
$ws.Range("E12").Value = '  +6.81%  '
$ws.Range("E13").Value = '  +12.30%  '
$ws.Range("D14").Value = '38.01'
$ws.Range("E14").Value = '  +7.91%  '
$ws.Range("D15").Value = '3.671.34'
$ws.Range("E15").Value = '  +2.18%  '
$ws.Range("D16").Value = '64.795.65'
$ws.Range("E16").Value = '  +0.62%  '
$ws.Range("D17").Value = '7.23'
$ws.Range("E17").Value = '  +6.85%  '
$ws.Range("D18").Value = '3.156.57'
$ws.Range("E18").Value = '  +2.14%  '
$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D19").Value = '0.111'
$ws.Range("E19").Value = '  +0.28%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").Value = '517.29'
$ws.Range("E20").Value = '  +6.65%  '
$ws.Range("D21").Value = '14.99'
$ws.Range("E21").Value = '  +6.96%  '
$ws.Range("D22").Value = '0.735'
$ws.Range("E22").Value = '  +8.74%  '
$ws.Range("D23").Value = '15.21'
$ws.Range("E23").Value = '  +6.82%  '
$ws.Range("E24").Value = '  +3.82%  '
$ws.Range("D25").Value = '85.30'
$ws.Range("E25").Value = '  +4.90%  '
$ws.Range("E26").Value = '  -0.16%  '
$ws.Range("D27").Value = '2.92'
$ws.Range("E27").Value = '  +4.26%  '
$ws.Range("D28").Value = '8.75'
$ws.Range("E28").Value = '  +8.87%  '
$ws.Range("E29").Value = '  +6.21%  '
$ws.Range("D30").Value = '27.94'
$ws.Range("E30").Value = '  +5.67%  '
$ws.Range("E31").Value = '  +0.11%  '
$ws.Range("D32").Value = '2.67'
$ws.Range("E32").Value = '  +7.79%  '
$ws.Range("E33").Value = '  +3.72%  '
$ws.Range("D34").Value = '6.10'
$ws.Range("E34").Value = '  +8.69%  '
$ws.Range("D35").Value = '6.57'
$ws.Range("E35").Value = '  +5.39%  '
$ws.Range("D36").Value = '55.70'
$ws.Range("E36").Value = '  +0.03%  '
$ws.Range("D37").Value = '486.13'
$ws.Range("E37").Value = '  +7.33%  '
$ws.Range("D38").Value = '0.0866'
$ws.Range("E38").Value = '  +5.92%  '
$ws.Range("D39").Value = '0.0422'
$ws.Range("E39").Value = '  +3.02%  '
$ws.Range("D40").Value = '2.96'
$ws.Range("D41").Value = '3.114.83'
$ws.Range("E41").Value = '  +5.06%  '
$ws.Range("D42").Value = '8.65'
$ws.Range("E42").Value = '  +4.98%  '
$ws.Range("E43").Value = '  +6.18%  '
$ws.Range("D44").Value = '0.294'
$ws.Range("E44").Value = '  +12.39%  '
$ws.Range("D45").Value = '2.43'
$ws.Range("E45").Value = '  +13.80%  '
$ws.Range("D46").Value = '29.24'
$ws.Range("E46").Value = '  +3.53%  '
$ws.Range("D47").Value = '0.0₃0576'
$ws.Range("E47").Value = '  +11.42%  '
$ws.Range("E48").Value = '  -0.03%  '
$ws.Range("E49").Value = '  +3.03%  '
$ws.Range("E50").Value = '  +9.89%  '
$ws.Range("D51").Value = '119.06'

# Restore the original (default) style on the cells we forced to text format.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = $origStyles[$addr]
}

